$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$value = 944138813511303040

for ($r = 5; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Formula = $value
}

$ws.Range("C17").Select()
